$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column 61 (BI). This shifts the existing
# "nom" column (BI) to BJ and the existing "url_produit" column (BJ) to BK,
# and also bumps the sheet dimension from BJ206 to BK206 automatically.
$ws.Columns.Item(61).Insert()

# New header cell for the freshly inserted column: a timestamp, matching
# the pattern used by every other column header in row 1.
$ws.Cells.Item(1, 61).Value = "2026-01-30 11:20:34"

# For every data row that already had a numeric price in column BH
# (rows 2-80), mirror that same price into the newly inserted column BI -
# this is the new "latest price" snapshot column lining up with BH.
# NOTE: use Value() (method-call syntax) to read the resolved scalar -
# plain property access returns the property descriptor, not the data.
for ($row = 2; $row -le 80; $row++) {
    $price = $ws.Cells.Item($row, 60).Value()
    $ws.Cells.Item($row, 61).Value = $price
}
